# Actualización automática 2025-10-31 14:30:09
# Apply updated sales figures across the three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-client figures by product group
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M4").Value  = 6099.16
$wsGrupo.Range("L5").Value  = 3847.3
$wsGrupo.Range("K13").Value = 832.3200000000001
$wsGrupo.Range("L13").Value = 1330.32
$wsGrupo.Range("M13").Value = 5374.16
$wsGrupo.Range("M16").Value = 1737.46
$wsGrupo.Range("M17").Value = -354.36
$wsGrupo.Range("M29").Value = 7379.76
$wsGrupo.Range("L37").Value = 1430.4
$wsGrupo.Range("K42").Value = 842.76
$wsGrupo.Range("K52").Value = 730.8

# Row 56 holds "count of non-zero-of-54" summary text per column
$wsGrupo.Range("K56").Value = "13 de 54"
$wsGrupo.Range("L56").Value = "7 de 54"
$wsGrupo.Range("M56").Value = "15 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" — monthly sales per client (October / col F)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value  = 9124.43
$wsMensual.Range("F5").Value  = 13542.31
$wsMensual.Range("F13").Value = 7536.8
$wsMensual.Range("F16").Value = 6465.16
$wsMensual.Range("F17").Value = -354.36
$wsMensual.Range("F29").Value = 12734.33
$wsMensual.Range("F37").Value = 7482.37
$wsMensual.Range("F42").Value = 929.16
$wsMensual.Range("F53").Value = 5829.37
$wsMensual.Range("F54").Value = 5829.37
$wsMensual.Range("F60").Value = 117403.54

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — budget compliance summary
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D10").Value = 13591.89
$wsCumpl.Range("E10").Value = -9710.810164656079
$wsCumpl.Range("F10").Value = 3.502089773114796

$wsCumpl.Range("D11").Value = 14723.57
$wsCumpl.Range("E11").Value = -2892.57
$wsCumpl.Range("F11").Value = 1.244490744653875

$wsCumpl.Range("D12").Value = 69797.50999999999
$wsCumpl.Range("E12").Value = -17134.38999999999
$wsCumpl.Range("F12").Value = 1.325358429200548

$wsCumpl.Range("D14").Value = 111310.11
$wsCumpl.Range("E14").Value = -12293.60338809385
$wsCumpl.Range("F14").Value = 1.124157110857066
